# Update "Sprint 1 - Bilan" sheet: mark the individual-grain management
# tasks (length, volume, pan, playback speed/pitch) and the ADSR/AS/AD
# envelope tasks as finished, log the time invested, and add the two new
# comments explaining the pitch/speed relationship and the early ADSR work.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1 - Bilan")

# Rows 9-12: grain length / volume / pan / loop (pitch) speed -> Finis, Court
$rows_court = @(9, 10, 11, 12, 18, 19)
foreach ($r in $rows_court) {
    $ws.Range("D$r").Value = "Finis"
    $ws.Range("E$r").Value = 1
    $ws.Range("F$r").Value = "Court"
}

# Row 12 gets an additional comment about pitch/speed
$ws.Range("G12").Value = "Vitesse de lecture = pitch dans ce context."

# Row 20: ADSR envelope -> Finis, Long, with a comment about being done early
$ws.Range("D20").Value = "Finis"
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = "Long"
$ws.Range("G20").Value = "ADSR a été fait plus tôt que prévu puisque façile a implémenter avec le MIDI"

# Update the active selection to match the saved view state
$ws.Range("F16").Select()
